# Auto-generated edit script updating leve profit calculations (H:N) across multiple sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 688.9
$ws.Range("I32").Value = 600.3333
$ws.Range("J32").Value = 704.5294
$ws.Range("K32").Value = 600.3333
$ws.Range("L32").Value = 704.5294
$ws.Range("M32").Value = -274.3333
$ws.Range("N32").Value = -1356.5294

$ws.Range("H113").Value = 2468.3333
$ws.Range("I113").Value = 2202.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2202.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1051.5
$ws.Range("N113").Value = -9508

$ws.Range("H118").Value = 1943.619
$ws.Range("I118").Value = 3097.8
$ws.Range("J118").Value = 1582.9375
$ws.Range("K118").Value = 9293.400000000001
$ws.Range("L118").Value = 4748.8125
$ws.Range("M118").Value = -7636.400000000001
$ws.Range("N118").Value = -8062.8125

$ws.Range("H125").Value = 974.75
$ws.Range("I125").Value = 507.5
$ws.Range("J125").Value = 1442
$ws.Range("K125").Value = 4567.5
$ws.Range("L125").Value = 12978
$ws.Range("M125").Value = -2107.5
$ws.Range("N125").Value = -17898

$ws.Range("H127").Value = 1766.0952
$ws.Range("I127").Value = 839.25
$ws.Range("J127").Value = 3001.889
$ws.Range("K127").Value = 2517.75
$ws.Range("L127").Value = 9005.667000000001
$ws.Range("M127").Value = 2442.25
$ws.Range("N127").Value = -18925.667

$ws.Range("H138").Value = 4363.146
$ws.Range("J138").Value = 4387.9
$ws.Range("L138").Value = 13163.7
$ws.Range("N138").Value = -23443.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15087.143
$ws.Range("I2").Value = 20702
$ws.Range("J2").Value = 1050
$ws.Range("K2").Value = 20702
$ws.Range("L2").Value = 1050
$ws.Range("M2").Value = -20589
$ws.Range("N2").Value = -1276

$ws.Range("H32").Value = 876158
$ws.Range("I32").Value = 1326832.4
$ws.Range("J32").Value = 15779.546
$ws.Range("K32").Value = 1326832.4
$ws.Range("L32").Value = 15779.546
$ws.Range("M32").Value = -1326545.4
$ws.Range("N32").Value = -16353.546

$ws.Range("H102").Value = 5813.5557
$ws.Range("I102").Value = 2460
$ws.Range("J102").Value = 10005.5
$ws.Range("K102").Value = 2460
$ws.Range("L102").Value = 10005.5
$ws.Range("M102").Value = -838
$ws.Range("N102").Value = -13249.5

$ws.Range("H116").Value = 15087.143
$ws.Range("I116").Value = 20702
$ws.Range("J116").Value = 1050
$ws.Range("K116").Value = 20702
$ws.Range("L116").Value = 1050
$ws.Range("M116").Value = -18408
$ws.Range("N116").Value = -5638

$ws.Range("H122").Value = 1347.9412
$ws.Range("I122").Value = 1132.5
$ws.Range("J122").Value = 1539.4445
$ws.Range("K122").Value = 3397.5
$ws.Range("L122").Value = 4618.333500000001
$ws.Range("M122").Value = -947.5
$ws.Range("N122").Value = -9518.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15087.143
$ws.Range("I3").Value = 20702
$ws.Range("J3").Value = 1050
$ws.Range("K3").Value = 20702
$ws.Range("L3").Value = 1050
$ws.Range("M3").Value = -20588
$ws.Range("N3").Value = -1278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1248.7826
$ws.Range("I31").Value = 1481.8
$ws.Range("J31").Value = 1069.5385
$ws.Range("K31").Value = 1481.8
$ws.Range("L31").Value = 1069.5385
$ws.Range("M31").Value = -1186.8
$ws.Range("N31").Value = -1659.5385

$ws.Range("H34").Value = 1248.7826
$ws.Range("I34").Value = 1481.8
$ws.Range("J34").Value = 1069.5385
$ws.Range("K34").Value = 1481.8
$ws.Range("L34").Value = 1069.5385
$ws.Range("M34").Value = -1279.8
$ws.Range("N34").Value = -1473.5385

$ws.Range("H94").Value = 674.2632
$ws.Range("I94").Value = 526.6667
$ws.Range("J94").Value = 742.38464
$ws.Range("K94").Value = 526.6667
$ws.Range("L94").Value = 742.38464
$ws.Range("M94").Value = -75.66669999999999
$ws.Range("N94").Value = -1644.38464

$ws.Range("H132").Value = 23813498
$ws.Range("I132").Value = 4437.3335
$ws.Range("K132").Value = 13312.0005
$ws.Range("M132").Value = -10782.0005

$ws.Range("H141").Value = 63065.2
$ws.Range("J141").Value = 28831.5
$ws.Range("L141").Value = 28831.5
$ws.Range("N141").Value = -39191.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4348.8
$ws.Range("I17").Value = 600
$ws.Range("J17").Value = 6848
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 20544
$ws.Range("M17").Value = -1631
$ws.Range("N17").Value = -20882

$ws.Range("H44").Value = 263.33334
$ws.Range("I44").Value = 263.33334
$ws.Range("K44").Value = 790.0000200000001
$ws.Range("M44").Value = -392.0000200000001

$ws.Range("H68").Value = 900.7308
$ws.Range("I68").Value = 736.5
$ws.Range("J68").Value = 950
$ws.Range("K68").Value = 2209.5
$ws.Range("L68").Value = 2850
$ws.Range("M68").Value = -1398.5
$ws.Range("N68").Value = -4472

$ws.Range("H71").Value = 900.7308
$ws.Range("I71").Value = 736.5
$ws.Range("J71").Value = 950
$ws.Range("K71").Value = 6628.5
$ws.Range("L71").Value = 8550
$ws.Range("M71").Value = -2572.5
$ws.Range("N71").Value = -16662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5820.8057
$ws.Range("I70").Value = 5850.6553
$ws.Range("J70").Value = 5697.143
$ws.Range("K70").Value = 5850.6553
$ws.Range("L70").Value = 5697.143
$ws.Range("M70").Value = -5580.6553
$ws.Range("N70").Value = -6237.143

$ws.Range("H73").Value = 5820.8057
$ws.Range("I73").Value = 5850.6553
$ws.Range("J73").Value = 5697.143
$ws.Range("K73").Value = 5850.6553
$ws.Range("L73").Value = 5697.143
$ws.Range("M73").Value = -4914.6553
$ws.Range("N73").Value = -7569.143

$ws.Range("H113").Value = 1877.3334
$ws.Range("I113").Value = 1562.75
$ws.Range("K113").Value = 1562.75
$ws.Range("M113").Value = 607.25

$ws.Range("H122").Value = 3595.45
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 3806.8125
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 11420.4375
$ws.Range("M122").Value = -5800
$ws.Range("N122").Value = -16320.4375

$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 4292
$ws.Range("I132").Value = 4376.6665
$ws.Range("J132").Value = 4139.6
$ws.Range("K132").Value = 13129.9995
$ws.Range("L132").Value = 12418.8
$ws.Range("M132").Value = -10599.9995
$ws.Range("N132").Value = -17478.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2628.55
$ws.Range("I100").Value = 1957.4166
$ws.Range("J100").Value = 3635.25
$ws.Range("K100").Value = 1957.4166
$ws.Range("L100").Value = 3635.25
$ws.Range("M100").Value = -1416.4166
$ws.Range("N100").Value = -4717.25

$ws.Range("H132").Value = 4084.44
$ws.Range("I132").Value = 3756.4443
$ws.Range("K132").Value = 11269.3329
$ws.Range("M132").Value = -8739.332900000001
